$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24
$ws.Cells.Item($row, 1).Value = 48
$ws.Cells.Item($row, 2).Value = "Title_12:05"
$ws.Cells.Item($row, 3).Value = "riya-morankar"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "edit1 to main"
$ws.Cells.Item($row, 6).NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = "2025-06-19"
